# Auto-generated edit script: applies numeric updates to the Typhon_Profits workbook
# per the commit "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1750
$ws.Range("J32").Value = 1750
$ws.Range("L32").Value = 1750
$ws.Range("N32").Value = -2402

$ws.Range("H97").Value = 2316.7778
$ws.Range("J97").Value = 2316.7778
$ws.Range("L97").Value = 6950.3334
$ws.Range("N97").Value = -7942.3334

$ws.Range("H106").Value = 1632.4
$ws.Range("I106").Value = 1019.7931
$ws.Range("J106").Value = 4593.3335
$ws.Range("K106").Value = 1019.7931
$ws.Range("L106").Value = 4593.3335
$ws.Range("M106").Value = -388.7931
$ws.Range("N106").Value = -5855.3335

$ws.Range("H132").Value = 6449.857
$ws.Range("I132").Value = 8476.223
$ws.Range("K132").Value = 25428.669
$ws.Range("M132").Value = -22898.669

$ws.Range("H137").Value = 46078.87
$ws.Range("I137").Value = 2940.2144
$ws.Range("J137").Value = 113183.445
$ws.Range("K137").Value = 8820.643199999999
$ws.Range("L137").Value = 339550.335
$ws.Range("M137").Value = -6270.643199999999
$ws.Range("N137").Value = -344650.335

$ws.Range("H141").Value = 3535
$ws.Range("I141").Value = 2879
$ws.Range("J141").Value = 4355
$ws.Range("K141").Value = 8637
$ws.Range("L141").Value = 13065
$ws.Range("M141").Value = -3457
$ws.Range("N141").Value = -23425

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1489.4286
$ws.Range("I2").Value = 1316.2858
$ws.Range("J2").Value = 2008.8572
$ws.Range("K2").Value = 1316.2858
$ws.Range("L2").Value = 2008.8572
$ws.Range("M2").Value = -1203.2858
$ws.Range("N2").Value = -2234.8572

$ws.Range("H116").Value = 1489.4286
$ws.Range("I116").Value = 1316.2858
$ws.Range("J116").Value = 2008.8572
$ws.Range("K116").Value = 1316.2858
$ws.Range("L116").Value = 2008.8572
$ws.Range("M116").Value = 977.7141999999999
$ws.Range("N116").Value = -6596.8572

$ws.Range("H119").Value = 25000
$ws.Range("J119").Value = 25000
$ws.Range("L119").Value = 25000
$ws.Range("N119").Value = -34676

$ws.Range("H122").Value = 2114
$ws.Range("I122").Value = 2114
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6342
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3892
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1489.4286
$ws.Range("I3").Value = 1316.2858
$ws.Range("J3").Value = 2008.8572
$ws.Range("K3").Value = 1316.2858
$ws.Range("L3").Value = 2008.8572
$ws.Range("M3").Value = -1202.2858
$ws.Range("N3").Value = -2236.8572

$ws.Range("H20").Value = 3558.2942
$ws.Range("I20").Value = 4040.3845
$ws.Range("K20").Value = 4040.3845
$ws.Range("M20").Value = -3793.3845

$ws.Range("H99").Value = 2020
$ws.Range("I99").Value = 1616.6666
$ws.Range("K99").Value = 1616.6666
$ws.Range("M99").Value = -118.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15300.852
$ws.Range("I31").Value = 37329.11
$ws.Range("J31").Value = 4286.722
$ws.Range("K31").Value = 37329.11
$ws.Range("L31").Value = 4286.722
$ws.Range("M31").Value = -37034.11
$ws.Range("N31").Value = -4876.722

$ws.Range("H34").Value = 15300.852
$ws.Range("I34").Value = 37329.11
$ws.Range("J34").Value = 4286.722
$ws.Range("K34").Value = 37329.11
$ws.Range("L34").Value = 4286.722
$ws.Range("M34").Value = -37127.11
$ws.Range("N34").Value = -4690.722

$ws.Range("H58").Value = 13203.829
$ws.Range("J58").Value = 39345.383
$ws.Range("L58").Value = 39345.383
$ws.Range("N58").Value = -39751.383

$ws.Range("H86").Value = 5383239.5
$ws.Range("I86").Value = 2379.4
$ws.Range("K86").Value = 2379.4
$ws.Range("M86").Value = -1256.4

$ws.Range("H89").Value = 5383239.5
$ws.Range("I89").Value = 2379.4
$ws.Range("K89").Value = 11897
$ws.Range("M89").Value = -6281

$ws.Range("H134").Value = 1385.0571
$ws.Range("I134").Value = 1058.6923
$ws.Range("J134").Value = 1577.909
$ws.Range("K134").Value = 3176.0769
$ws.Range("L134").Value = 4733.727000000001
$ws.Range("M134").Value = -641.0769
$ws.Range("N134").Value = -9803.727000000001

$ws.Range("H136").Value = 13203.829
$ws.Range("J136").Value = 39345.383
$ws.Range("L136").Value = 118036.149
$ws.Range("N136").Value = -123136.149

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 250000500
$ws.Range("J86").Value = 500000000
$ws.Range("L86").Value = 1500000000
$ws.Range("N86").Value = -1500002372

$ws.Range("H89").Value = 250000500
$ws.Range("J89").Value = 500000000
$ws.Range("L89").Value = 4500000000
$ws.Range("N89").Value = -4500011856

$ws.Range("H98").Value = 1027.2858
$ws.Range("J98").Value = 738.8
$ws.Range("L98").Value = 2216.4
$ws.Range("N98").Value = -5212.4

$ws.Range("H131").Value = 772.53
$ws.Range("J131").Value = 779.7292
$ws.Range("L131").Value = 2339.1876
$ws.Range("N131").Value = -12419.1876

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 121213360
$ws.Range("I122").Value = 55556840
$ws.Range("J122").Value = 200001170
$ws.Range("K122").Value = 166670520
$ws.Range("L122").Value = 600003510
$ws.Range("M122").Value = -166668070
$ws.Range("N122").Value = -600008410

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5738.3076
$ws.Range("I61").Value = 3227.8572
$ws.Range("J61").Value = 8667.166999999999
$ws.Range("K61").Value = 3227.8572
$ws.Range("L61").Value = 8667.166999999999
$ws.Range("M61").Value = -3025.8572
$ws.Range("N61").Value = -9071.166999999999

$ws.Range("H93").Value = 2554.4443
$ws.Range("I93").Value = 2723.75
$ws.Range("J93").Value = 1200
$ws.Range("K93").Value = 2723.75
$ws.Range("L93").Value = 1200
$ws.Range("M93").Value = -1475.75
$ws.Range("N93").Value = -3696

$ws.Range("H113").Value = 5738.3076
$ws.Range("I113").Value = 3227.8572
$ws.Range("J113").Value = 8667.166999999999
$ws.Range("K113").Value = 3227.8572
$ws.Range("L113").Value = 8667.166999999999
$ws.Range("M113").Value = -1057.8572
$ws.Range("N113").Value = -13007.167

$ws.Range("H122").Value = 1964395.2
$ws.Range("I122").Value = 2453781
$ws.Range("J122").Value = 6852.5
$ws.Range("K122").Value = 7361343
$ws.Range("L122").Value = 20557.5
$ws.Range("M122").Value = -7358893
$ws.Range("N122").Value = -25457.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4638.6665
$ws.Range("I62").Value = 3539
$ws.Range("J62").Value = 5424.143
$ws.Range("K62").Value = 3539
$ws.Range("L62").Value = 5424.143
$ws.Range("M62").Value = -2915
$ws.Range("N62").Value = -6672.143

$ws.Range("H65").Value = 4638.6665
$ws.Range("I65").Value = 3539
$ws.Range("J65").Value = 5424.143
$ws.Range("K65").Value = 17695
$ws.Range("L65").Value = 27120.715
$ws.Range("M65").Value = -14575
$ws.Range("N65").Value = -33360.715

$ws.Range("H122").Value = 1867.6
$ws.Range("I122").Value = 1908.4117
$ws.Range("J122").Value = 1780.875
$ws.Range("K122").Value = 5725.2351
$ws.Range("L122").Value = 5342.625
$ws.Range("M122").Value = -3275.2351
$ws.Range("N122").Value = -10242.625
